$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.213.36'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").Value = '1.856.78'
$ws.Range("E3").Value = '  -2.33%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4799'
$ws.Range("E7").Value = '  -2.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2801'
$ws.Range("E8").Value = '  -4.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06457'
$ws.Range("E9").Value = '  -3.37%  '
$ws.Range("D10").Value = '1.862.33'
$ws.Range("E10").Value = '  -2.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07384'
$ws.Range("E11").Value = '  +0.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.25'
$ws.Range("E12").Value = '  -4.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.032'
$ws.Range("E13").Value = '  -2.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '87.08'
$ws.Range("E14").Value = '  -1.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6451'
$ws.Range("E15").Value = '  -3.48%  '
$ws.Range("D16").Value = '30.181.19'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9999'
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.14'
$ws.Range("E18").Value = '  -2.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007557'
$ws.Range("E19").Value = '  -4.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '223.47'
$ws.Range("E20").Value = '  +11.84%  '
$ws.Range("D21").Value = '2.100.62'
$ws.Range("E21").Value = '  -2.12%  '
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.281'
$ws.Range("E23").Value = '  -2.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.075'
$ws.Range("E24").Value = '  -1.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.195'
$ws.Range("E25").Value = '  -3.27%  '
$ws.Range("E26").Value = '  +0.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.49'
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.924'
$ws.Range("E28").Value = '  -1.01%  '
$ws.Range("E29").Value = '  -3.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09201'
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.238'
$ws.Range("E31").Value = '  -2.44%  '
$ws.Range("E32").Value = '  -4.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04964'
$ws.Range("E33").Value = '  -3.78%  '
$ws.Range("E34").Value = '  +3.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7237'
$ws.Range("E35").Value = '  -2.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01830'
$ws.Range("E37").Value = '  -0.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.597'
$ws.Range("E38").Value = '  -3.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8991'
$ws.Range("E39").Value = '  -2.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.038'
$ws.Range("E40").Value = '  -1.53%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.894'
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '105.92'
$ws.Range("E42").Value = '  -0.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4248'
$ws.Range("E43").Value = '  -3.73%  '
$ws.Range("E44").Value = '  +0.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1302'
$ws.Range("E45").Value = '  -5.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.264'
$ws.Range("E46").Value = '  -4.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '63.61'
$ws.Range("E47").Value = '  -8.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.494'
$ws.Range("E48").Value = '  +6.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.716'
$ws.Range("E49").Value = '  -3.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.83'
$ws.Range("E50").Value = '  -3.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05634'
$ws.Range("E51").Value = '  -3.55%  '
